$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Activate()

# Insert a new column before column N (14) - this shifts N:P -> O:Q
$ws.Columns("N:N").Insert()

# Set the width of the newly inserted column N to match column M (stored width 11)
$ws.Columns("N:N").ColumnWidth = 10.17

# Update selection to K13 as seen in the diff
$ws.Range("K13").Select()
